$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9355671405792236
$ws.Range("B1").Value = 1.085562944412231
$ws.Range("C1").Value = 1.397915840148926
$ws.Range("D1").Value = 3.071463346481323
$ws.Range("E1").Value = 4.291028022766113
